# Apply predictor-label edits described by the diff:
# several "(per capita)" / bare predictor names get wrapped with ln(...)
# and the "Livestock AB Consumption (kg per capita)" label gets a
# (typo'd) bracket swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new predictor text (column C, "Predictor")
$updates = @{
    "C2"  = "ln(GDP [dollars per capita])"
    "C5"  = "ln(Tourism - Inbound [per capita])"
    "C6"  = "ln(Migrant Population [per capita])"
    "C9"  = "ln(ProMed Mentions [per capita])"
    "C10" = "ln(AB Exports [dollars per capita])"
    "C11" = "ln(Publication Bias Index [per capita])"
    "C12" = "Livestock AB Consumption [kg per capita)"
    "C13" = "ln(ProMed Mentions [per capita])"
    "C15" = "ln(Publication Bias Index [per capita])"
    "C16" = "ln(Population)"
    "C17" = "ln(GDP [dollars per capita])"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
